$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old "Sugerir 2 Features (...)" entries in column A (rows 2-6)
$ws.Range("A2:A6").ClearContents()

# Write the new task list into column D (rows 2-5)
$ws.Range("D2").Value = "Debugs"
$ws.Range("D3").Value = "Correct bugs"
$ws.Range("D4").Value = "Analyse how to implement features"
$ws.Range("D5").Value = "Make java program to encapsulate tasks"

# Update the selection to match the new layout
$ws.Range("D6").Select()
